$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '37.476.17'
$ws.Range('E2').Value = '  +5.26%  '
Set-TextValue 'D3' '2.055.88'
$ws.Range('E3').Value = '  +3.77%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '252.84'
$ws.Range('E5').Value = '  +2.98%  '
Set-TextValue 'D6' '0.653'
$ws.Range('E6').Value = '  +2.75%  '
Set-TextValue 'D7' '65.30'
$ws.Range('E7').Value = '  +13.34%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +6.60%  '
Set-TextValue 'D10' '59.76'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('E11').Value = '  +4.99%  '
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('E13').Value = '  -2.45%  '
Set-TextValue 'D14' '14.90'
$ws.Range('E14').Value = '  +2.89%  '
Set-TextValue 'D15' '22.77'
$ws.Range('E15').Value = '  +26.42%  '
$ws.Range('E16').Value = '  +3.91%  '
$ws.Range('E17').Value = '  +6.08%  '
Set-TextValue 'D18' '2.051.17'
$ws.Range('E18').Value = '  +3.57%  '
Set-TextValue 'D19' '37.369.25'
$ws.Range('E19').Value = '  +5.14%  '
Set-TextValue 'D20' '73.65'
$ws.Range('E20').Value = '  +3.15%  '
Set-TextValue 'D21' '0.0₃0878'
$ws.Range('E21').Value = '  +4.36%  '
Set-TextValue 'D22' '5.49'
$ws.Range('E22').Value = '  +6.28%  '
Set-TextValue 'D23' '240.09'
$ws.Range('E23').Value = '  +3.06%  '
Set-TextValue 'D24' '2.66'
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +4.83%  '
Set-TextValue 'D27' '10.17'
$ws.Range('E27').Value = '  +11.92%  '
Set-TextValue 'D28' '162.48'
$ws.Range('E28').Value = '  -1.30%  '
$ws.Range('E29').Value = '  +4.42%  '
Set-TextValue 'D30' '0.120'
$ws.Range('E30').Value = '  +25.02%  '
Set-TextValue 'D31' '5.30'
$ws.Range('E31').Value = '  +8.93%  '
$ws.Range('E32').Value = '  +3.14%  '
Set-TextValue 'D33' '1.22'
$ws.Range('E33').Value = '  +8.59%  '
$ws.Range('E34').Value = '  +9.00%  '
$ws.Range('E35').Value = '  +6.01%  '
Set-TextValue 'D36' '2.45'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  +4.35%  '
$ws.Range('E39').Value = '  +15.67%  '
Set-TextValue 'D40' '2.99'
$ws.Range('E40').Value = '  +33.23%  '
$ws.Range('E41').Value = '  +17.74%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D42' '3.03'
$ws.Range('E42').Value = '  +4.53%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '1.25'
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('E44').Value = '  +6.10%  '
Set-TextValue 'D45' '17.34'
$ws.Range('E45').Value = '  +8.59%  '
$ws.Range('E46').Value = '  +3.59%  '
Set-TextValue 'D47' '97.08'
$ws.Range('E47').Value = '  +6.02%  '
$ws.Range('E48').Value = '  +5.44%  '
Set-TextValue 'D49' '1.425.15'
$ws.Range('E49').Value = '  +4.07%  '
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D51' '3.74'
$ws.Range('E51').Value = '  +4.89%  '
